$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: replace the visible text of a table cell without touching the
# trailing paragraph/cell-end mark, and without letting Find.Execute's
# document-wide search scope touch other (possibly duplicate-text) cells.
function Set-CellText($row, $col, $newText) {
    $cellRange = $t.Cell($row, $col).Range
    $innerRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $innerRange.Text = $newText
}

# --- Header row text updates (row 1) ---
# Col 2: "チャイの売上合計 (単位)" -> "チャイの売上合計 (ユニット数)"
Set-CellText 1 2 "チャイの売上合計 (ユニット数)"

# Col 3: "職人チャイ販売 (ユニット)" -> "Artisanal Chai の売上 (ユニット数)"
Set-CellText 1 3 "Artisanal Chai の売上 (ユニット数)"

# Col 4: "事前に作成されたチャイの売上 (単位)" -> "既製チャイの売上 (ユニット数)"
Set-CellText 1 4 "既製チャイの売上 (ユニット数)"

# Col 5: make the header run bold, and update its text to add "(ビュー)"
$t.Cell(1, 5).Range.Font.Bold = 1
Set-CellText 1 5 "ソーシャル メディア エンゲージメント (ビュー)"

# --- Date column fixes (column 1) ---
# Row 7 (May): "3/31/2023" -> "2023/5/31"
Set-CellText 7 1 "2023/5/31"

# Row 9 (July): "3/30/2023" -> "2023/7/30"
Set-CellText 9 1 "2023/7/30"

# Row 11 (September): "2023 年 9 月 5 日" -> "2023/9/30"
Set-CellText 11 1 "2023/9/30"

# Row 13 (November): "2020/11/30" -> "2023/11/30"
Set-CellText 13 1 "2023/11/30"
